# "Upload tree DBH and hollow data mk2"
#
# The "Large?" (col E) / "Tree Class" (col C) / notes (col L) columns for
# the second transect (rows 40-77) used a lowercase "yes" answer that was
# a separate, inconsistent entry in the shared-string table from the
# "Yes" answer used everywhere else in the sheet (e.g. column E rows 2-19).
# This pass normalizes every lowercase "yes" cell to the canonical "Yes"
# so the workbook only has one spelling for that answer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell that currently holds the lowercase "yes" value (rows 40-77,
# columns C/E/L). Listed explicitly so we touch exactly the cells that
# need normalizing and nothing else.
$yesCells = @(
    "C40", "E40", "L40",
    "C41", "E41", "L41",
    "E42",
    "C43", "E43", "L43",
    "E44",
    "C45", "E45", "L45",
    "C46", "E46", "L46",
    "C47", "E47", "L47",
    "C48", "E48", "L48",
    "C49", "E49", "L49",
    "C50", "E50", "L50",
    "C51", "E51", "L51",
    "C52", "E52", "L52",
    "C53", "E53", "L53",
    "C54", "E54", "L54",
    "C55", "E55", "L55",
    "C56", "E56",
    "C57", "E57", "L57",
    "C58", "L58",
    "C59", "L59",
    "C60", "L60",
    "C61", "L61",
    "C62", "L62",
    "C63", "L63",
    "C64", "L64",
    "C65", "L65",
    "C66", "L66",
    "C67", "L67",
    "C68", "L68",
    "C69", "L69",
    "C70", "L70",
    "C71", "L71",
    "C72", "L72",
    "C73", "L73",
    "C74", "L74",
    "C75", "L75",
    "C76", "L76",
    "C77", "L77"
)

foreach ($addr in $yesCells) {
    $ws.Range($addr).Value2 = "Yes"
}

# Move the saved selection/scroll position to where the editor left off.
$ws.Range("N75").Select()
